# Update "paises.xlsx" (sheet "Pais") with the latest COVID snapshot data.
# Source refresh: 16 de Octubre de 2020, 04:16 -> 05:33.
#
# The sheet lists one country per row, sorted by "Casos totales" (col B)
# descending. Column A holds a fixed rank label (stable per row), while
# B-H hold the per-country counters that get refreshed on every pull.
# Because B is updated, a handful of neighbouring countries swap ranks
# (their A-label moves up/down a row or two) even though most rows keep
# the same counters-column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 05:33"

# --- Row 26 (rank 30, Pakistan) — refreshed counters ---
$ws.Range("B26").Value = 321877
$ws.Range("C26").Value = 659
$ws.Range("D26").Value = 305835
$ws.Range("E26").Value = 9421
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = 6621

# --- Rows 30-31 swap: Canada and Belgica trade places ---
$ws.Range("A30").Value = "Belgica"
$ws.Range("B30").Value = 191959
$ws.Range("C30").Value = 10448
$ws.Range("D30").Value = 20720
$ws.Range("E30").Value = 160912
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 49
$ws.Range("H30").Value = 10327

$ws.Range("A31").Value = "Canada"
$ws.Range("B31").Value = 191732
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 161490
$ws.Range("E31").Value = 20543
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 9699

# --- Row 44 (rank 48, Kazajistan) — refreshed counters ---
$ws.Range("B44").Value = 109202
$ws.Range("C44").Value = 108
$ws.Range("D44").Value = 104801
$ws.Range("E44").Value = 2633
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 1768

# --- Rows 53-55 cascade: Honduras moves up ahead of Bielorrusia & China ---
$ws.Range("A53").Value = "Honduras"
$ws.Range("B53").Value = 86089
$ws.Range("C53").Value = 631
$ws.Range("D53").Value = 34099
$ws.Range("E53").Value = 49438
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 19
$ws.Range("H53").Value = 2552

$ws.Range("A54").Value = "Bielorrusia"
$ws.Range("B54").Value = 85734
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 78583
$ws.Range("E54").Value = 6235
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 916

$ws.Range("A55").Value = "China"
$ws.Range("B55").Value = 85646
$ws.Range("C55").Value = 24
$ws.Range("D55").Value = 80759
$ws.Range("E55").Value = 253
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 4634

# --- Row 56 (rank 60, Venezuela) — refreshed counters, same rank ---
$ws.Range("B56").Value = 85469
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 77689
$ws.Range("E56").Value = 7060
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 720

# --- Row 85 (rank 89, Australia) — refreshed counters ---
$ws.Range("B85").Value = 27371
$ws.Range("C85").Value = 9
$ws.Range("D85").Value = 25062
$ws.Range("E85").Value = 1405
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 904

# --- Row 153 (rank 157, Belice) — refreshed counters ---
$ws.Range("B153").Value = 2682
$ws.Range("C153").Value = 63
$ws.Range("D153").Value = 1612
$ws.Range("E153").Value = 1029
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 41
